# Yearly coverage in scenario 3b
#
# The "Platform Coverage" sheet's first coverage row (row 2, the yearly MDA
# coverage for the 5-15 age band) only had values on even years
# (2018, 2020, 2022, ...). Fill in the missing odd years too so every year
# from 2018-2040 (columns H:AD) carries the same 0.6 coverage value.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Platform Coverage")

$ws1.Range("H2:AD2").Value = 0.6

# Make "Platform Coverage" the active/selected sheet (it was "MarketShare"
# before), and leave the selection on the last filled-in cell, AC2, with the
# view scrolled so column Q is the left-most visible column.
$ws1.Activate()
$ws1.Range("AC2").Select()
$excel.ActiveWindow.ScrollColumn = $ws1.Range("Q1").Column

$wb.Save()
